# Update cryptocurrency price/volume data per the Feb 16 2023 symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'322.57"
$ws.Range("E2").Value = "'8.53%"
$ws.Range("D3").Value = "'50.40"
$ws.Range("E3").Value = "'20.53%"
$ws.Range("D4").Value = "'5.345"
$ws.Range("E4").Value = "'6.54%"
$ws.Range("D5").Value = "'0.08164"
$ws.Range("E5").Value = "'8.53%"
$ws.Range("D6").Value = "'4.602"
$ws.Range("E6").Value = "'5.27%"
$ws.Range("D7").Value = "'1.674"
$ws.Range("E7").Value = "'5.51%"
$ws.Range("D8").Value = "'1.206"
$ws.Range("E8").Value = "'30.61%"
$ws.Range("D9").Value = "'0.1340"
$ws.Range("E9").Value = "'12.47%"
$ws.Range("E10").Value = "'7.09%"
$ws.Range("D11").Value = "'0.09678"
$ws.Range("E11").Value = "'8.13%"
$ws.Range("D12").Value = "'0.04473"
$ws.Range("E12").Value = "'9.30%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("D14").Value = "'0.001318"
$ws.Range("E14").Value = "'3.18%"
$ws.Range("D15").Value = "'0.005892"
$ws.Range("E15").Value = "'0.33%"
$ws.Range("D16").Value = "'3.390"
$ws.Range("E16").Value = "'1.42%"
$ws.Range("D17").Value = "'2.435"
$ws.Range("E17").Value = "'1.40%"
$ws.Range("E18").Value = "'2.40%"
$ws.Range("D19").Value = "'8.123"
$ws.Range("E19").Value = "'0.46%"
$ws.Range("D20").Value = "'0.1418"
$ws.Range("E20").Value = "'2.73%"
$ws.Range("D21").Value = "'0.3050"
$ws.Range("E21").Value = "'-1.69%"
$ws.Range("D22").Value = "'0.04310"
$ws.Range("E22").Value = "'5.78%"
$ws.Range("D23").Value = "'0.001304"
$ws.Range("E23").Value = "'3.02%"
$ws.Range("E24").Value = "'9.57%"
$ws.Range("D25").Value = "'0.0001349"
$ws.Range("E25").Value = "'9.56%"
$ws.Range("D26").Value = "'0.0003536"
$ws.Range("E26").Value = "'-5.05%"
$ws.Range("D38").Value = "'0.02765"
$ws.Range("E38").Value = "'14.58%"
$ws.Range("D39").Value = "'0.05595"
$ws.Range("E39").Value = "'7.32%"
$ws.Range("D40").Value = "'0.006294"
$ws.Range("E40").Value = "'-0.16%"
$ws.Range("D41").Value = "'0.007684"
$ws.Range("E41").Value = "'-1.60%"
$ws.Range("D42").Value = "'0.1449"
$ws.Range("E42").Value = "'9.07%"
$ws.Range("D43").Value = "'0.007675"
$ws.Range("E43").Value = "'3.79%"
$ws.Range("D44").Value = "'0.008069"
$ws.Range("E44").Value = "'3.80%"
$ws.Range("D45").Value = "'0.3511"
$ws.Range("E45").Value = "'18.37%"
$ws.Range("D46").Value = "'0.00006817"
$ws.Range("E46").Value = "'3.47%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.16%"
$ws.Range("E48").Value = "'93.31%"
$ws.Range("D49").Value = "'0.003996"
$ws.Range("E49").Value = "'-4.91%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.16%"
